# Actualización automática 2025-10-06 15:30:20
# Updates the "CUMPLIMIENTO MENSUAL" sheet (3rd sheet) of the workbook:
#  - refreshes VENTA (D) figures for several GRUPO rows and recomputes
#    POR CUMPLIR (E) and CUMPLIMIENTO (F)
#  - removes the "SAL SOLUBLE" detail row (old row 14)
#  - recomputes the TOTAL row (now row 14 after the deletion)
#  - narrows/widens a few columns

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

function Set-Cumplimiento($row, $presupuesto, $venta) {
    $porCumplir = $presupuesto - $venta
    if ($presupuesto -ne 0) {
        $cumplimiento = $venta / $presupuesto
    } else {
        $cumplimiento = 0
    }
    $ws.Cells.Item($row, 4).Value2 = $venta
    $ws.Cells.Item($row, 5).Value2 = $porCumplir
    $ws.Cells.Item($row, 6).Value2 = $cumplimiento
}

# Row 3 - 240X80 PORCELANATO
Set-Cumplimiento 3 8834.57354940916 190.08

# Row 4 - FREGADEROS DE COCINA
Set-Cumplimiento 4 521.61144263264 0

# Row 6 - INODOROS
Set-Cumplimiento 6 814.123430808873 0

# Row 8 - NO RESURTIBLES
Set-Cumplimiento 8 480.217743214072 0

# Row 11 - PIEDRA SINTERIZADA
Set-Cumplimiento 11 2922.22458185274 358.23

# Row 12 - PORCELANATO (PRESUPUESTO also changes here)
$ws.Cells.Item(12, 3).Value2 = 21701.27
Set-Cumplimiento 12 21701.27 377.45

# Remove the "SAL SOLUBLE" row entirely (old row 14); TOTAL (old row 15)
# shifts up to become the new row 14.
$ws.Rows.Item(14).Delete()

# Recompute the TOTAL row (now row 14) from the remaining detail rows (2-13)
$totalPresupuesto = 0
$totalVenta = 0
for ($r = 2; $r -le 13; $r++) {
    $totalPresupuesto += $ws.Cells.Item($r, 3).Value2
    $totalVenta += $ws.Cells.Item($r, 4).Value2
}
$ws.Cells.Item(14, 3).Value2 = $totalPresupuesto
Set-Cumplimiento 14 $totalPresupuesto $totalVenta

# Column width adjustments (D, E, F). ColumnWidth includes Excel's standard
# padding of ~0.8333 characters over the stored OOXML width, so subtract it
# to land on the exact target widths.
$pad = 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth = 12 - $pad
$ws.Columns.Item(5).ColumnWidth = 22 - $pad
$ws.Columns.Item(6).ColumnWidth = 25 - $pad
